$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-8
# from 2023-11-03 (45233) to 2023-11-13 (45243)
for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 3).Value = 45243
}
